$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '31.201.38'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +2.62%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.998.22'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +6.73%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7816'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +65.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '257.50'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +4.96%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9987'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3527'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +22.90%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '29.12'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +33.16%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07041'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +8.43%  '
$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.8617'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +18.37%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08209'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +5.24%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.999.41'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +6.82%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '100.98'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.63%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.592'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +8.16%  '
$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.56'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +18.77%  '
$ws.Range('B17').Value = 'BitcoinCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '274.13'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -3.48%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '31.210.90'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.72%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.961'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +11.89%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007957'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +6.35%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.258.55'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +6.88%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9992'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9997'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.145'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +12.81%  '
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.09'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +11.55%  '
$ws.Range('B26').Value = 'Stellar'
$ws.Range('C26').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1489'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +53.97%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '164.21'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.18%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.97'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +5.35%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.377'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +25.40%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.610'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +8.03%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.627'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +9.40%  '
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.358'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +2.85%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.441'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +6.95%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05219'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +8.65%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.231'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +9.42%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7761'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +12.57%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.812'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +3.46%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02015'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +6.19%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.909'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.60%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.715'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +6.65%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '79.77'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +4.21%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4728'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +12.02%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.156'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +10.21%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '106.40'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +5.45%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8526'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.37%  '
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9987'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.739'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +10.29%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.999'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.69%  '
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4338'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +10.81%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.84'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +5.31%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1195'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +13.56%  '
